$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking price values so Excel keeps them as text
# (loop individually - a comma-joined multi-area Range only applies
#  NumberFormat to the first area in this runtime)
foreach ($addr in @('D4','D5','D6','D7','D8','D9','D10','D11','D13','D14','D16','D18','D19','D22','D25','D26','D27','D28','D29','D30','D36','D38','D39','D45','D46','D47','D48','D50')) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range('D2').Value = '30.621.41'
$ws.Range('E2').Value = '  +2.28%  '

# Row 3
$ws.Range('D3').Value = '1.675.08'
$ws.Range('E3').Value = '  +2.69%  '

# Row 4
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.03%  '

# Row 5
$ws.Range('D5').Value = '219.92'
$ws.Range('E5').Value = '  +2.60%  '

# Row 6
$ws.Range('D6').Value = '0.532'
$ws.Range('E6').Value = '  +2.81%  '

# Row 7
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.05%  '

# Row 8
$ws.Range('D8').Value = '29.63'
$ws.Range('E8').Value = '  +3.88%  '

# Row 9
$ws.Range('D9').Value = '0.264'
$ws.Range('E9').Value = '  +2.58%  '

# Row 10
$ws.Range('D10').Value = '0.0651'
$ws.Range('E10').Value = '  +7.00%  '

# Row 11
$ws.Range('D11').Value = '0.0904'
$ws.Range('E11').Value = '  -0.20%  '

# Row 12
$ws.Range('D12').Value = '1.916.23'
$ws.Range('E12').Value = '  +2.82%  '

# Row 13
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').Value = '0.613'
$ws.Range('E13').Value = '  +9.14%  '

# Row 14
$ws.Range('D14').Value = '10.17'
$ws.Range('E14').Value = '  +11.09%  '

# Row 15
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '1.648.12'
$ws.Range('E15').Value = '  +1.45%  '

# Row 16
$ws.Range('D16').Value = '4.03'
$ws.Range('E16').Value = '  +4.97%  '

# Row 17
$ws.Range('D17').Value = '30.647.24'
$ws.Range('E17').Value = '  +2.33%  '

# Row 18
$ws.Range('D18').Value = '66.37'
$ws.Range('E18').Value = '  +3.58%  '

# Row 19
$ws.Range('D19').Value = '243.25'
$ws.Range('E19').Value = '  +0.60%  '

# Row 20
$ws.Range('D20').Value = '0.0₃0727'
$ws.Range('E20').Value = '  +3.83%  '

# Row 21
$ws.Range('E21').Value = '  -0.15%  '

# Row 22
$ws.Range('D22').Value = '4.26'
$ws.Range('E22').Value = '  +3.23%  '

# Row 23
$ws.Range('E23').Value = '  +2.72%  '

# Row 24
$ws.Range('E24').Value = '  -0.31%  '

# Row 25
$ws.Range('D25').Value = '158.57'
$ws.Range('E25').Value = '  +0.10%  '

# Row 26
$ws.Range('B26').Value = 'Stellar'
$ws.Range('C26').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D26').Value = '0.113'
$ws.Range('E26').Value = '  +3.11%  '

# Row 27
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '15.87'
$ws.Range('E27').Value = '  +2.41%  '

# Row 28
$ws.Range('D28').Value = '6.69'
$ws.Range('E28').Value = '  +1.62%  '

# Row 29
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.05%  '

# Row 30
$ws.Range('D30').Value = '0.0495'
$ws.Range('E30').Value = '  +1.83%  '

# Row 31
$ws.Range('E31').Value = '  +4.11%  '

# Row 32
$ws.Range('E32').Value = '  +2.62%  '

# Row 33
$ws.Range('E33').Value = '  +3.92%  '

# Row 34
$ws.Range('D34').Value = '1.491.32'
$ws.Range('E34').Value = '  +4.61%  '

# Row 35
$ws.Range('E35').Value = '  +7.53%  '

# Row 36
$ws.Range('D36').Value = '84.34'
$ws.Range('E36').Value = '  +11.42%  '

# Row 37
$ws.Range('E37').Value = '  -0.59%  '

# Row 38
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.0178'
$ws.Range('E38').Value = '  +5.73%  '

# Row 39
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').Value = '0.596'
$ws.Range('E39').Value = '  +8.03%  '

# Row 40
$ws.Range('E40').Value = '  -3.36%  '

# Row 41
$ws.Range('E41').Value = '  -0.14%  '

# Row 42
$ws.Range('E42').Value = '  +1.63%  '

# Row 43
$ws.Range('E43').Value = '  +2.08%  '

# Row 44
$ws.Range('E44').Value = '  -0.77%  '

# Row 45
$ws.Range('D45').Value = '1.02'
$ws.Range('E45').Value = '  +0.38%  '

# Row 46
$ws.Range('D46').Value = '0.999'
$ws.Range('E46').Value = '  -0.03%  '

# Row 47
$ws.Range('D47').Value = '51.70'
$ws.Range('E47').Value = '  -0.69%  '

# Row 48
$ws.Range('D48').Value = '5.53'
$ws.Range('E48').Value = '  +3.58%  '

# Row 49
$ws.Range('D49').Value = '1.808.72'
$ws.Range('E49').Value = '  +2.14%  '

# Row 50
$ws.Range('D50').Value = '94.63'
$ws.Range('E50').Value = '  +5.09%  '

# Row 51
$ws.Range('D51').Value = '0.0₆0116'
$ws.Range('E51').Value = '  +1.89%  '

